# Weekly data refresh: a new price observation is inserted as the first
# record of the "Espinaca" (Arica) series, pushing every existing record
# (previously rows 53-82) down by one row (now rows 54-83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 53 - everything currently at/after row 53
# shifts down one row (row 82 -> 83), matching the new sheet dimension
# A1:R83.
$ws.Rows("53:53").Insert()

# Populate the new row 53 with the latest weekly observation. The
# descriptive columns (market, region, product, quality grade, unit,
# origin, classification, etc.) repeat the same values used throughout
# this block of rows; only the date and the price/volume figures are new.
$ws.Range("A53").Value = 1
$ws.Range("B53").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C53").Value = "Arica y Parinacota"
$ws.Range("D53").Value = 44873
$ws.Range("E53").Value = 15
$ws.Range("F53").Value = 100112012
$ws.Range("G53").Value = "Espinaca"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 250
$ws.Range("K53").Value = 1500
$ws.Range("L53").Value = 1700
$ws.Range("M53").Value = 1600
$ws.Range("N53").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O53").Value = "Región de Arica y Parinacota"
$ws.Range("P53").Value = 533
$ws.Range("Q53").Value = 3
$ws.Range("R53").Value = "Hortaliza"
